# Natmi following Dr Hou advice
# Update Pdpn-Clec1b LR-pair data: full 3x3 sending/target cluster combinations
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pdpn"
$ws.Range("C2").Value = "Clec1b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 2.141780333333334
$ws.Range("H2").Value = 6.425341
$ws.Range("I2").Value = 0.02869229445178705
$ws.Range("J2").Value = 0.02869229445178705
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.82600633333333
$ws.Range("N2").Value = 38.478019
$ws.Range("O2").Value = 0.8780027849272567
$ws.Range("P2").Value = 0.8780027849272567
$ws.Range("Q2").Value = 27.47048811994211
$ws.Range("R2").Value = 247.234393079479
$ws.Range("S2").Value = 0.02519191443462191
$ws.Range("T2").Value = 0.0251919144346219

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pdpn"
$ws.Range("C3").Value = "Clec1b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 2.141780333333334
$ws.Range("H3").Value = 6.425341
$ws.Range("I3").Value = 0.02869229445178705
$ws.Range("J3").Value = 0.02869229445178705
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.702433
$ws.Range("N3").Value = 5.107299
$ws.Range("O3").Value = 0.1165398547533384
$ws.Range("P3").Value = 0.1165398547533384
$ws.Range("Q3").Value = 3.646237518217667
$ws.Range("R3").Value = 32.816137663959
$ws.Range("S3").Value = 0.003343795827951279
$ws.Range("T3").Value = 0.003343795827951279

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pdpn"
$ws.Range("C4").Value = "Clec1b"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 2.141780333333334
$ws.Range("H4").Value = 6.425341
$ws.Range("I4").Value = 0.02869229445178705
$ws.Range("J4").Value = 0.02869229445178705
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.079722
$ws.Range("N4").Value = 0.239166
$ws.Range("O4").Value = 0.005457360319405017
$ws.Range("P4").Value = 0.005457360319405017
$ws.Range("Q4").Value = 0.170747011734
$ws.Range("R4").Value = 1.536723105606
$ws.Range("S4").Value = 0.0001565841892138674
$ws.Range("T4").Value = 0.0001565841892138674

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pdpn"
$ws.Range("C5").Value = "Clec1b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 45.835794
$ws.Range("H5").Value = 137.507382
$ws.Range("I5").Value = 0.6140378064974858
$ws.Range("J5").Value = 0.6140378064974858
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.82600633333333
$ws.Range("N5").Value = 38.478019
$ws.Range("O5").Value = 0.8780027849272567
$ws.Range("P5").Value = 0.8780027849272567
$ws.Range("Q5").Value = 587.890184137362
$ws.Range("R5").Value = 5291.011657236258
$ws.Range("S5").Value = 0.5391269041554165
$ws.Range("T5").Value = 0.5391269041554165

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pdpn"
$ws.Range("C6").Value = "Clec1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 45.835794
$ws.Range("H6").Value = 137.507382
$ws.Range("I6").Value = 0.6140378064974858
$ws.Range("J6").Value = 0.6140378064974858
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.702433
$ws.Range("N6").Value = 5.107299
$ws.Range("O6").Value = 0.1165398547533384
$ws.Range("P6").Value = 0.1165398547533384
$ws.Range("Q6").Value = 78.032368286802
$ws.Range("R6").Value = 702.2913145812181
$ws.Range("S6").Value = 0.07155987678227549
$ws.Range("T6").Value = 0.07155987678227549

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pdpn"
$ws.Range("C7").Value = "Clec1b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 45.835794
$ws.Range("H7").Value = 137.507382
$ws.Range("I7").Value = 0.6140378064974858
$ws.Range("J7").Value = 0.6140378064974858
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.079722
$ws.Range("N7").Value = 0.239166
$ws.Range("O7").Value = 0.005457360319405017
$ws.Range("P7").Value = 0.005457360319405017
$ws.Range("Q7").Value = 3.654121169268
$ws.Range("R7").Value = 32.887090523412
$ws.Range("S7").Value = 0.003351025559793875
$ws.Range("T7").Value = 0.003351025559793875

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Pdpn"
$ws.Range("C8").Value = "Clec1b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 26.66895966666667
$ws.Range("H8").Value = 80.006879
$ws.Range("I8").Value = 0.357269899050727
$ws.Range("J8").Value = 0.357269899050727
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.82600633333333
$ws.Range("N8").Value = 38.478019
$ws.Range("O8").Value = 0.8780027849272567
$ws.Range("P8").Value = 0.8780027849272567
$ws.Range("Q8").Value = 342.0562455880778
$ws.Range("R8").Value = 3078.5062102927
$ws.Range("S8").Value = 0.3136839663372182
$ws.Range("T8").Value = 0.3136839663372182

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Pdpn"
$ws.Range("C9").Value = "Clec1b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 26.66895966666667
$ws.Range("H9").Value = 80.006879
$ws.Range("I9").Value = 0.357269899050727
$ws.Range("J9").Value = 0.357269899050727
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.702433
$ws.Range("N9").Value = 5.107299
$ws.Range("O9").Value = 0.1165398547533384
$ws.Range("P9").Value = 0.1165398547533384
$ws.Range("Q9").Value = 45.40211701220233
$ws.Range("R9").Value = 408.619053109821
$ws.Range("S9").Value = 0.04163618214311159
$ws.Range("T9").Value = 0.04163618214311159

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Pdpn"
$ws.Range("C10").Value = "Clec1b"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 26.66895966666667
$ws.Range("H10").Value = 80.006879
$ws.Range("I10").Value = 0.357269899050727
$ws.Range("J10").Value = 0.357269899050727
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.079722
$ws.Range("N10").Value = 0.239166
$ws.Range("O10").Value = 0.005457360319405017
$ws.Range("P10").Value = 0.005457360319405017
$ws.Range("Q10").Value = 2.126102802546
$ws.Range("R10").Value = 19.134925222914
$ws.Range("S10").Value = 0.001949750570397274
$ws.Range("T10").Value = 0.001949750570397274

